# cambios en selection sort y pantalla
#
# Adds a new "Promedios" (averages) block (rows 31-44) comparing sort
# algorithm timings for array sizes 10 and 5, and updates the active
# selection / scroll position on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers for the "size Array 10" comparison table -----------------
# Write in the exact order the strings were first introduced so the
# generated shared-string table matches (Quick sort time, Sort Algorithm for
# size Array 10 , Tree Sort time, Selection Sort time, Sort Algorithm for
# size Array 5, Bogo Sort, Promedios).
$ws.Range("C32").Value = "Quick sort time"
$ws.Range("C31").Value = "Sort Algorithm for size Array 10 "
$ws.Range("D32").Value = "Tree Sort time"
$ws.Range("E32").Value = "Selection Sort time"
$ws.Range("G31").Value = "Sort Algorithm for size Array 5"
$ws.Range("G32").Value = "Bogo Sort"

# --- Timing samples (rows 33-42) -------------------------------------------
$data = @(
  @(1.11, 2.02, 9.03, 14.89),
  @(1.21, 2.02, 9.04, 5.04),
  @(1.21, 2.02, 9.04, 28.98),
  @(1.51, 2.02, 9.04, 26.56),
  @(1.31, 2.01, 9.04, 50.93),
  @(1.21, 2.02, 9.04, 2.22),
  @(1.41, 2.02, 9.04, 4.83),
  @(1.41, 2.02, 9.04, 3.62),
  @(1.21, 2.02, 9.03, 2.62),
  @(1.21, 2.02, 9.04, 28.59)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = 33 + $i
  $ws.Cells.Item($row, 3).Value = $data[$i][0]
  $ws.Cells.Item($row, 4).Value = $data[$i][1]
  $ws.Cells.Item($row, 5).Value = $data[$i][2]
  $ws.Cells.Item($row, 7).Value = $data[$i][3]
}

# --- Underlined, otherwise empty marker cell -------------------------------
$ws.Range("I43").Font.Underline = 2

# --- Averages row -----------------------------------------------------------
$ws.Range("B44").Value = "Promedios"
$ws.Range("C44").Formula = "=AVERAGE(C33:C42)"
$ws.Range("D44").Formula = "=AVERAGE(D33:D42)"
$ws.Range("E44").Formula = "=AVERAGE(E33:E42)"
$ws.Range("G44").Formula = "=AVERAGE(G33:G42)"

# --- Selection / scroll position -------------------------------------------
$ws.Range("I43").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 2
